$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16, shifting rows 16:22 down to 17:23
$ws.Rows.Item(16).Insert()

# Fill in the new row 16 with values (a new weekly record)
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value = 45176
$ws.Cells.Item(16, 4).Style = $ws.Cells.Item(17, 4).Style
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(17, 4).NumberFormat
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 100112017
$ws.Cells.Item(16, 7).Value = "Ramas de apio"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 250
$ws.Cells.Item(16, 11).Value = 5000
$ws.Cells.Item(16, 12).Value = 6000
$ws.Cells.Item(16, 13).Value = 5500
$ws.Cells.Item(16, 14).Value = "`$/atado 7 kilos"
$ws.Cells.Item(16, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(16, 16).Value = 5500
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(16, 18).Value = "Hortaliza"
